$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 577
$ws.Range("J97").Value = 577
$ws.Range("L97").Value = 1731
$ws.Range("N97").Value = -2723
$ws.Range("H113").Value = 5512.5
$ws.Range("H129").Value = 2248
$ws.Range("I129").Value = 1500
$ws.Range("K129").Value = 4500
$ws.Range("M129").Value = 500
$ws.Range("H132").Value = 8507.826999999999
$ws.Range("I132").Value = 8471.370000000001
$ws.Range("K132").Value = 25414.11
$ws.Range("M132").Value = -22884.11
$ws.Range("H135").Value = 641.3570999999999
$ws.Range("I135").Value = 666.0769
$ws.Range("K135").Value = 5994.6921
$ws.Range("M135").Value = -3459.6921
$ws.Range("H137").Value = 2029.75
$ws.Range("I137").Value = 965.5
$ws.Range("J137").Value = 3094
$ws.Range("K137").Value = 2896.5
$ws.Range("L137").Value = 9282
$ws.Range("M137").Value = -346.5
$ws.Range("N137").Value = -14382
$ws.Range("H138").Value = 2115.25
$ws.Range("I138").Value = 1024.1
$ws.Range("J138").Value = 7571
$ws.Range("K138").Value = 3072.3
$ws.Range("L138").Value = 22713
$ws.Range("M138").Value = 2067.7
$ws.Range("N138").Value = -32993
$ws.Range("H141").Value = 1147.3334
$ws.Range("I141").Value = 1147.3334
$ws.Range("K141").Value = 3442.0002
$ws.Range("M141").Value = 1737.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1950.625
$ws.Range("I61").Value = 1353.8182
$ws.Range("K61").Value = 1353.8182
$ws.Range("M61").Value = -1141.8182
$ws.Range("H74").Value = 2071.1304
$ws.Range("I74").Value = 1437.8823
$ws.Range("J74").Value = 3865.3333
$ws.Range("K74").Value = 1437.8823
$ws.Range("L74").Value = 3865.3333
$ws.Range("M74").Value = -563.8823
$ws.Range("N74").Value = -5613.3333
$ws.Range("H77").Value = 2071.1304
$ws.Range("I77").Value = 1437.8823
$ws.Range("J77").Value = 3865.3333
$ws.Range("K77").Value = 7189.4115
$ws.Range("L77").Value = 19326.6665
$ws.Range("M77").Value = -2821.4115
$ws.Range("N77").Value = -28062.6665
$ws.Range("H119").Value = 19000
$ws.Range("J119").Value = 19000
$ws.Range("L119").Value = 19000
$ws.Range("N119").Value = -28676
$ws.Range("H132").Value = 1696.4667
$ws.Range("I132").Value = 1696.4667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5089.4001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2559.4001
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1950.625
$ws.Range("I136").Value = 1353.8182
$ws.Range("K136").Value = 4061.4546
$ws.Range("M136").Value = -1511.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 3333.3333
$ws.Range("I75").Value = 3333.3333
$ws.Range("K75").Value = 3333.3333
$ws.Range("M75").Value = -2397.3333
$ws.Range("H78").Value = 3333.3333
$ws.Range("I78").Value = 3333.3333
$ws.Range("K78").Value = 9999.999899999999
$ws.Range("M78").Value = -5319.999899999999
$ws.Range("H80").Value = 693.0833
$ws.Range("I80").Value = 334.8889
$ws.Range("J80").Value = 1767.6666
$ws.Range("K80").Value = 334.8889
$ws.Range("L80").Value = 1767.6666
$ws.Range("M80").Value = 663.1111000000001
$ws.Range("N80").Value = -3763.6666
$ws.Range("H83").Value = 693.0833
$ws.Range("I83").Value = 334.8889
$ws.Range("J83").Value = 1767.6666
$ws.Range("K83").Value = 1674.4445
$ws.Range("L83").Value = 8838.333000000001
$ws.Range("M83").Value = 3317.5555
$ws.Range("N83").Value = -18822.333
$ws.Range("H94").Value = 897.7273
$ws.Range("I94").Value = 887.5
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 887.5
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -436.5
$ws.Range("N94").Value = -1902
$ws.Range("H134").Value = 1952.7
$ws.Range("I134").Value = 793.2857
$ws.Range("J134").Value = 4658
$ws.Range("K134").Value = 2379.8571
$ws.Range("L134").Value = 13974
$ws.Range("M134").Value = 155.1428999999998
$ws.Range("N134").Value = -19044
$ws.Range("H140").Value = 88260
$ws.Range("J140").Value = 88260
$ws.Range("L140").Value = 88260
$ws.Range("N140").Value = -98620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6944.5264
$ws.Range("I31").Value = 2353
$ws.Range("K31").Value = 2353
$ws.Range("M31").Value = -2058
$ws.Range("H34").Value = 6944.5264
$ws.Range("I34").Value = 2353
$ws.Range("K34").Value = 2353
$ws.Range("M34").Value = -2151
$ws.Range("H58").Value = 2645.1765
$ws.Range("I58").Value = 1351.8182
$ws.Range("J58").Value = 5016.3335
$ws.Range("K58").Value = 1351.8182
$ws.Range("L58").Value = 5016.3335
$ws.Range("M58").Value = -1148.8182
$ws.Range("N58").Value = -5422.3335
$ws.Range("H94").Value = 11004.667
$ws.Range("J94").Value = 14507
$ws.Range("L94").Value = 14507
$ws.Range("N94").Value = -15409
$ws.Range("H122").Value = 1065
$ws.Range("I122").Value = 1078
$ws.Range("K122").Value = 3234
$ws.Range("M122").Value = -784
$ws.Range("H132").Value = 1764.25
$ws.Range("I132").Value = 1838.6923
$ws.Range("K132").Value = 5516.0769
$ws.Range("M132").Value = -2986.0769
$ws.Range("H134").Value = 3814.7778
$ws.Range("I134").Value = 2222
$ws.Range("J134").Value = 7000.3335
$ws.Range("K134").Value = 6666
$ws.Range("L134").Value = 21001.0005
$ws.Range("M134").Value = -4131
$ws.Range("N134").Value = -26071.0005
$ws.Range("H136").Value = 2645.1765
$ws.Range("I136").Value = 1351.8182
$ws.Range("J136").Value = 5016.3335
$ws.Range("K136").Value = 4055.4546
$ws.Range("L136").Value = 15049.0005
$ws.Range("M136").Value = -1505.4546
$ws.Range("N136").Value = -20149.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83.3
$ws.Range("I2").Value = 77.789474
$ws.Range("J2").Value = 188
$ws.Range("K2").Value = 466.736844
$ws.Range("L2").Value = 1128
$ws.Range("M2").Value = -353.736844
$ws.Range("N2").Value = -1354
$ws.Range("H57").Value = 1589.2
$ws.Range("I57").Value = 1299
$ws.Range("K57").Value = 3897
$ws.Range("M57").Value = -3338
$ws.Range("H59").Value = 443.33334
$ws.Range("I59").Value = 443.33334
$ws.Range("K59").Value = 1330.00002
$ws.Range("M59").Value = -790.0000199999999
$ws.Range("H60").Value = 1736
$ws.Range("I60").Value = 227.5
$ws.Range("K60").Value = 682.5
$ws.Range("M60").Value = -431.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1498.2667
$ws.Range("I102").Value = 1194.7084
$ws.Range("K102").Value = 1194.7084
$ws.Range("M102").Value = 427.2916
$ws.Range("H107").Value = 421.54544
$ws.Range("I107").Value = 93.85714
$ws.Range("J107").Value = 995
$ws.Range("K107").Value = 93.85714
$ws.Range("L107").Value = 995
$ws.Range("M107").Value = 1826.14286
$ws.Range("N107").Value = -4835

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8204.625
$ws.Range("J7").Value = 9449.5
$ws.Range("L7").Value = 9449.5
$ws.Range("N7").Value = -9673.5
$ws.Range("H16").Value = 1200
$ws.Range("I16").Value = 1083.3334
$ws.Range("K16").Value = 1083.3334
$ws.Range("M16").Value = -913.3334
$ws.Range("H40").Value = 6247.5
$ws.Range("I40").Value = 4995.5
$ws.Range("K40").Value = 4995.5
$ws.Range("M40").Value = -4859.5
$ws.Range("H126").Value = 8204.625
$ws.Range("J126").Value = 9449.5
$ws.Range("L126").Value = 28348.5
$ws.Range("N126").Value = -33288.5
$ws.Range("H132").Value = 3302.0908
$ws.Range("I132").Value = 3144.8
$ws.Range("K132").Value = 9434.400000000001
$ws.Range("M132").Value = -6904.400000000001
$ws.Range("H136").Value = 1498.5
$ws.Range("I136").Value = 1498.5
$ws.Range("K136").Value = 4495.5
$ws.Range("M136").Value = -1945.5
$ws.Range("H141").Value = 37857.5
$ws.Range("J141").Value = 37857.5
$ws.Range("L141").Value = 37857.5
$ws.Range("N141").Value = -48217.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 480.7
$ws.Range("I107").Value = 480.7
$ws.Range("K107").Value = 1442.1
$ws.Range("M107").Value = 477.9000000000001
$ws.Range("H120").Value = 38000
$ws.Range("J120").Value = 38000
$ws.Range("L120").Value = 38000
$ws.Range("N120").Value = -47676
$ws.Range("H132").Value = 1411.2759
$ws.Range("I132").Value = 1318.8214
$ws.Range("K132").Value = 3956.4642
$ws.Range("M132").Value = -1426.4642
$ws.Range("H136").Value = 3579.8484
$ws.Range("I136").Value = 2833.818
$ws.Range("K136").Value = 8501.454000000002
$ws.Range("M136").Value = -5951.454000000002
$ws.Range("H141").Value = 195903.33
$ws.Range("J141").Value = 135084
$ws.Range("L141").Value = 135084
$ws.Range("N141").Value = -145444
